$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Group" column (I) ---
# Header cell: same text/style as the other header cells
$ws.Range("I1").Value = "Group"
$ws.Range("H1").Copy()
[void]$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data: first data row belongs to group "A"; second data row has no group
$ws.Range("I2").Value = "A"

# --- Re-stamp the existing data cells (A2:H3) ---
# (mirrors the style bookkeeping churn seen when the sheet is regenerated
# after adding the required "Group" field)
$ws.Range("A2:H3").Style = "Normal"

# --- Row / selection bookkeeping ---
$ws.Rows.Item(1).RowHeight = 13.8
[void]$ws.Range("I3").Select()

Write-Host "done"
